# The stock price table contained yearly rows starting at EOG-FY-1994.
# Remove the earliest 8 fiscal years (EOG-FY-1994 through EOG-FY-2001),
# which are rows 2 through 9 of the sheet, shifting the remaining data
# (EOG-FY-2002 .. EOG-FY-2021) up to start at row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:9").Delete()
